$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3329.3333
$ws.Range("I86").Value = 1577.5
$ws.Range("K86").Value = 1577.5
$ws.Range("M86").Value = -454.5
$ws.Range("H87").Value = 97600
$ws.Range("J87").Value = 97600
$ws.Range("L87").Value = 97600
$ws.Range("N87").Value = -100096
$ws.Range("H89").Value = 3329.3333
$ws.Range("I89").Value = 1577.5
$ws.Range("K89").Value = 7887.5
$ws.Range("M89").Value = -2271.5
$ws.Range("H90").Value = 97600
$ws.Range("J90").Value = 97600
$ws.Range("L90").Value = 292800
$ws.Range("N90").Value = -305280
$ws.Range("H99").Value = 470
$ws.Range("I99").Value = 342.75
$ws.Range("J99").Value = 724.5
$ws.Range("K99").Value = 1028.25
$ws.Range("L99").Value = 2173.5
$ws.Range("M99").Value = 469.75
$ws.Range("N99").Value = -5169.5
$ws.Range("H112").Value = 2472.432
$ws.Range("J112").Value = 2472.432
$ws.Range("L112").Value = 7417.295999999999
$ws.Range("N112").Value = -9633.295999999998
$ws.Range("H118").Value = 407030
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H129").Value = 1981.5
$ws.Range("I129").Value = 1468
$ws.Range("K129").Value = 4404
$ws.Range("M129").Value = 596
$ws.Range("H137").Value = 13576.311
$ws.Range("I137").Value = 3463.875
$ws.Range("J137").Value = 17428.666
$ws.Range("K137").Value = 10391.625
$ws.Range("L137").Value = 52285.99800000001
$ws.Range("M137").Value = -7841.625
$ws.Range("N137").Value = -57385.99800000001
$ws.Range("H138").Value = 2571.2903
$ws.Range("I138").Value = 1553.8704
$ws.Range("J138").Value = 3980.0256
$ws.Range("K138").Value = 4661.6112
$ws.Range("L138").Value = 11940.0768
$ws.Range("M138").Value = 478.3887999999997
$ws.Range("N138").Value = -22220.0768

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 219.25
$ws.Range("I5").Value = 283.83334
$ws.Range("K5").Value = 283.83334
$ws.Range("M5").Value = -171.83334
$ws.Range("H132").Value = 3979.9495
$ws.Range("I132").Value = 2984.3428
$ws.Range("K132").Value = 8953.028399999999
$ws.Range("M132").Value = -6423.028399999999

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 219.25
$ws.Range("I4").Value = 283.83334
$ws.Range("K4").Value = 283.83334
$ws.Range("M4").Value = -168.83334
$ws.Range("H20").Value = 3597.8333
$ws.Range("J20").Value = 2679.7144
$ws.Range("L20").Value = 2679.7144
$ws.Range("N20").Value = -3173.7144
$ws.Range("H128").Value = 4000
$ws.Range("I128").Value = 4000
$ws.Range("K128").Value = 12000
$ws.Range("M128").Value = -9510
$ws.Range("H134").Value = 9939.261
$ws.Range("I134").Value = 3519.7273
$ws.Range("K134").Value = 10559.1819
$ws.Range("M134").Value = -8024.1819

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H31").Value = 43483336
$ws.Range("I31").Value = 76925896
$ws.Range("J31").Value = 7999.2
$ws.Range("K31").Value = 76925896
$ws.Range("L31").Value = 7999.2
$ws.Range("M31").Value = -76925601
$ws.Range("N31").Value = -8589.200000000001
$ws.Range("H34").Value = 43483336
$ws.Range("I34").Value = 76925896
$ws.Range("J34").Value = 7999.2
$ws.Range("K34").Value = 76925896
$ws.Range("L34").Value = 7999.2
$ws.Range("M34").Value = -76925694
$ws.Range("N34").Value = -8403.200000000001
$ws.Range("H58").Value = 8263.825999999999
$ws.Range("I58").Value = 9432.615
$ws.Range("J58").Value = 6744.4
$ws.Range("K58").Value = 9432.615
$ws.Range("L58").Value = 6744.4
$ws.Range("M58").Value = -9229.615
$ws.Range("N58").Value = -7150.4
$ws.Range("H132").Value = 3255.7966
$ws.Range("I132").Value = 2052.6667
$ws.Range("K132").Value = 6158.000100000001
$ws.Range("M132").Value = -3628.000100000001
$ws.Range("H134").Value = 3516.25
$ws.Range("I134").Value = 3242.8306
$ws.Range("J134").Value = 6742.6
$ws.Range("K134").Value = 9728.4918
$ws.Range("L134").Value = 20227.8
$ws.Range("M134").Value = -7193.4918
$ws.Range("N134").Value = -25297.8
$ws.Range("H136").Value = 8263.825999999999
$ws.Range("I136").Value = 9432.615
$ws.Range("J136").Value = 6744.4
$ws.Range("K136").Value = 28297.845
$ws.Range("L136").Value = 20233.2
$ws.Range("M136").Value = -25747.845
$ws.Range("N136").Value = -25333.2

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 264.2
$ws.Range("J23").Value = 279.66666
$ws.Range("L23").Value = 838.9999799999999
$ws.Range("N23").Value = -1308.99998
$ws.Range("H81").Value = 1000
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 1000
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H86").Value = 428.57144
$ws.Range("I86").Value = 200
$ws.Range("K86").Value = 600
$ws.Range("M86").Value = 586
$ws.Range("H89").Value = 428.57144
$ws.Range("I89").Value = 200
$ws.Range("K89").Value = 1800
$ws.Range("M89").Value = 4128
$ws.Range("H107").Value = 2320.158
$ws.Range("J107").Value = 2617.6875
$ws.Range("L107").Value = 7853.0625
$ws.Range("N107").Value = -11693.0625
$ws.Range("H132").Value = 1955.2222
$ws.Range("J132").Value = 2257.7144
$ws.Range("L132").Value = 20319.4296
$ws.Range("N132").Value = -25379.4296

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4332.683
$ws.Range("I132").Value = 2194.862
$ws.Range("J132").Value = 9499.083000000001
$ws.Range("K132").Value = 6584.586
$ws.Range("L132").Value = 28497.249
$ws.Range("M132").Value = -4054.586
$ws.Range("N132").Value = -33557.249

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 6000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 6000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 6000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -6226
$ws.Range("H7").Value = 56574.95
$ws.Range("I7").Value = 70994.60000000001
$ws.Range("K7").Value = 70994.60000000001
$ws.Range("M7").Value = -70882.60000000001
$ws.Range("H16").Value = 3698.8
$ws.Range("I16").Value = 4248.5
$ws.Range("K16").Value = 4248.5
$ws.Range("M16").Value = -4078.5
$ws.Range("H22").Value = 3088.6487
$ws.Range("I22").Value = 1975.8
$ws.Range("J22").Value = 4397.8823
$ws.Range("K22").Value = 1975.8
$ws.Range("L22").Value = 4397.8823
$ws.Range("M22").Value = -1680.8
$ws.Range("N22").Value = -4987.8823
$ws.Range("H27").Value = 3088.6487
$ws.Range("I27").Value = 1975.8
$ws.Range("J27").Value = 4397.8823
$ws.Range("K27").Value = 1975.8
$ws.Range("L27").Value = 4397.8823
$ws.Range("M27").Value = -1868.8
$ws.Range("N27").Value = -4611.8823
$ws.Range("H28").Value = 6000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 6000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 6000
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -6464
$ws.Range("H37").Value = 6000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 6000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 6000
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -6214
$ws.Range("H46").Value = 5176.0557
$ws.Range("I46").Value = 1750.3334
$ws.Range("J46").Value = 5487.485
$ws.Range("K46").Value = 1750.3334
$ws.Range("L46").Value = 5487.485
$ws.Range("M46").Value = -1562.3334
$ws.Range("N46").Value = -5863.485
$ws.Range("H55").Value = 395.26666
$ws.Range("I55").Value = 173
$ws.Range("J55").Value = 649.2857
$ws.Range("K55").Value = 173
$ws.Range("L55").Value = 649.2857
$ws.Range("M55").Value = 0
$ws.Range("N55").Value = -995.2857
$ws.Range("H109").Value = 79000
$ws.Range("J109").Value = 79000
$ws.Range("L109").Value = 79000
$ws.Range("N109").Value = -81774
$ws.Range("H126").Value = 56574.95
$ws.Range("I126").Value = 70994.60000000001
$ws.Range("K126").Value = 212983.8
$ws.Range("M126").Value = -210513.8
$ws.Range("H136").Value = 2907.7273
$ws.Range("I136").Value = 2535.3704
$ws.Range("J136").Value = 4583.3335
$ws.Range("K136").Value = 7606.111199999999
$ws.Range("L136").Value = 13750.0005
$ws.Range("M136").Value = -5056.111199999999
$ws.Range("N136").Value = -18850.0005

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 6684.7
$ws.Range("I14").Value = 8504
$ws.Range("J14").Value = 6229.875
$ws.Range("K14").Value = 8504
$ws.Range("L14").Value = 6229.875
$ws.Range("M14").Value = -8336
$ws.Range("N14").Value = -6565.875
$ws.Range("H26").Value = 17999.5
$ws.Range("I26").Value = 17999.5
$ws.Range("K26").Value = 17999.5
$ws.Range("M26").Value = -17706.5
$ws.Range("H113").Value = 1020
$ws.Range("I113").Value = 830.25
$ws.Range("J113").Value = 1399.5
$ws.Range("K113").Value = 2490.75
$ws.Range("L113").Value = 4198.5
$ws.Range("M113").Value = -320.75
$ws.Range("N113").Value = -8538.5
$ws.Range("H132").Value = 1962.6875
$ws.Range("I132").Value = 1288.8413
$ws.Range("K132").Value = 3866.5239
$ws.Range("M132").Value = -1336.5239
$ws.Range("H136").Value = 2014.079
$ws.Range("I136").Value = 616.75
$ws.Range("K136").Value = 1850.25
$ws.Range("M136").Value = 699.75
